# Trade #38 closed at 2026-02-17 12:39:19 - unknown UNKNOWN +0.000%
# Updates summary metrics, strategy status, and appends the new closed
# trade row to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.74
$wsSummary.Range("B4").Value = 0.73
$wsSummary.Range("B6").Value = 38
$wsSummary.Range("B7").Value = 15
$wsSummary.Range("B9").Value = 39.47

# ---------------------------------------------------------------------
# Sheet: Strategy Status (row 4 = MarketMaking)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.74
$wsStatus.Range("D4").Value = 38
$wsStatus.Range("E4").Value = 0.73
$wsStatus.Range("F4").Value = 0.74
$wsStatus.Range("G4").Value = 39.47

# ---------------------------------------------------------------------
# New trade row data (Trade #38)
# Note: the Date column (B) holds a plain text value (e.g. "2026-02-17")
# in the source file, not a real date. A leading apostrophe forces Excel
# to keep it as text instead of auto-converting it to a date serial.
# ---------------------------------------------------------------------
$tradeRow = 39
$values = @(38, "'2026-02-17", "12:39:13", "MarketMaking", "UP", 0.95, 0.97, "CLOSED", 2.1053, 0.02, 100.74, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

function Add-TradeRow($ws, $row, $vals) {
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
    $ws.Cells.Item($row, 10).Value = $vals[9]
    $ws.Cells.Item($row, 11).Value = $vals[10]
    $ws.Cells.Item($row, 12).Value = $vals[11]
    $ws.Cells.Item($row, 13).Value = $vals[12]
    $ws.Cells.Item($row, 14).Value = $vals[13]
    $ws.Cells.Item($row, 15).Value = $vals[14]
    $ws.Cells.Item($row, 16).Value = $vals[15]
    $ws.Cells.Item($row, 17).Value = $vals[16]
}

# ---------------------------------------------------------------------
# Sheet: All Trades
# ---------------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades $tradeRow $values

# ---------------------------------------------------------------------
# Sheet: MarketMaking
# ---------------------------------------------------------------------
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking $tradeRow $values

Write-Output "Applied trade #38 update"
